$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header area updates ---
$ws.Range("E11").Value = 624000
$ws.Range("F13").Value = 6

# --- Swap "Novedad de Ingreso" / "Novedad de Retiro" headers (H15 / I15) ---
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# --- Rebuild the worker statement-of-account rows ---
# Remove two rows from the top of the existing data block (rows 16-17). This
# shifts the remaining rows up by two, so the special "closing" row style
# (borders) that used to sit on row 29 now lands correctly on row 27, the new
# last row of the (now 12-row) table, and rows 34/35 below shift up to 32/33.
$ws.Rows("16:17").Delete()

$doc1 = "84453334"
$name1 = "RONALD ENRIQUE CURIEL DE LA HOZ"
$doc2 = "71777922"
$name2 = "JUAN CARLOS GARCIA GONZALEZ"

$periods = @("2505", "2504", "2503", "2502", "2501", "2412")

$workers = @(
    @{ Doc = $doc1; Name = $name1 },
    @{ Doc = $doc2; Name = $name2 }
)

$r = 16
foreach ($worker in $workers) {
    foreach ($p in $periods) {
        $ws.Range("B" + $r).Value = "CC"

        $ws.Range("C" + $r).NumberFormat = "@"
        $ws.Range("C" + $r).Value = $worker.Doc

        $ws.Range("D" + $r).Value = $worker.Name

        $ws.Range("E" + $r).NumberFormat = "@"
        $ws.Range("E" + $r).Value = $p

        $ws.Range("F" + $r).Value = 52000
        $ws.Range("G" + $r).Value = 1300000

        $r++
    }
}
